# Availability grid: "bulk select" for days/slots re-jiggered which Day column
# each existing booking landed in (column width of "Monday" vs "Thursday" swapped
# too, since the grid widened the day that now carries the moved bookings).
#
# Excel's COM ColumnWidth is expressed in "characters" of the Normal-style font
# and is internally offset by 5/6 (0.8333...) of a character from the stored
# OOXML <col width> value, so we subtract that constant to land exactly on the
# target stored widths of 12 and 37.
$wb = $excel.ActiveWorkbook

$kine = "kine_1 | agu | kine:javi | s1 | n=1"
$fono = "fono_1 | agu | fono:maca | s2 | n=1"
$widthOffset = 0.8333333333333334

# --- s1 (kine schedule) ---------------------------------------------------
$ws = $wb.Worksheets.Item("s1")
$ws.Columns.Item(2).ColumnWidth = 12 - $widthOffset   # Monday column: 37 -> 12
$ws.Columns.Item(5).ColumnWidth = 37 - $widthOffset   # Thursday column: 12 -> 37

$ws.Range("C4").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("D10").ClearContents()

$ws.Range("D5").Value = $kine
$ws.Range("C10").Value = $kine
$ws.Range("E10").Value = $kine

# --- s2 (fono schedule) ----------------------------------------------------
$ws = $wb.Worksheets.Item("s2")
$ws.Columns.Item(2).ColumnWidth = 12 - $widthOffset   # Monday column: 37 -> 12
$ws.Columns.Item(5).ColumnWidth = 37 - $widthOffset   # Thursday column: 12 -> 37

$ws.Range("C2").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("D7").ClearContents()

$ws.Range("D6").Value = $fono
$ws.Range("C8").Value = $fono
$ws.Range("E8").Value = $fono

# --- Therapists summary sheet ----------------------------------------------
$ws = $wb.Worksheets.Item("Therapists")
$ws.Range("D7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("D25").ClearContents()
$ws.Range("C28").ClearContents()

$ws.Range("D17").Value = $fono
$ws.Range("C19").Value = $kine
$ws.Range("C23").Value = $kine
$ws.Range("D24").Value = $fono
$ws.Range("D35").Value = $fono
$ws.Range("C37").Value = $kine

# --- Patients summary sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Patients")
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("C28").ClearContents()

$ws.Range("C17").Value = $fono
$ws.Range("C19").Value = $kine
$ws.Range("C23").Value = $kine
$ws.Range("C24").Value = $fono
$ws.Range("C35").Value = $fono
$ws.Range("C37").Value = $kine
